$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($CellRef, $Value)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "44.372.16"
$ws.Range("E2").Value = "  +5.27%  "
$ws.Range("D3").Value = "2.283.01"
$ws.Range("E3").Value = "  +1.69%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue "D5" "230.96"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("E6").Value = "  -0.44%  "
Set-TextValue "D7" "60.84"
$ws.Range("E7").Value = "  -1.52%  "
$ws.Range("E8").Value = "  -0.05%  "
Set-TextValue "D9" "0.424"
$ws.Range("E9").Value = "  +4.74%  "
Set-TextValue "D10" "0.0942"
$ws.Range("E10").Value = "  +4.88%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "2.623.33"
$ws.Range("E12").Value = "  +1.73%  "
Set-TextValue "D13" "24.24"
$ws.Range("E13").Value = "  +9.00%  "
Set-TextValue "D14" "15.57"
$ws.Range("E14").Value = "  -1.19%  "
Set-TextValue "D15" "5.83"
$ws.Range("E15").Value = "  +3.79%  "
Set-TextValue "D16" "0.808"
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "2.280.31"
$ws.Range("E17").Value = "  +1.69%  "
$ws.Range("D18").Value = "44.222.62"
$ws.Range("E18").Value = "  +5.00%  "
$ws.Range("D19").Value = "0.0₃0940"
$ws.Range("E19").Value = "  +3.75%  "
Set-TextValue "D20" "73.26"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("E21").Value = "  +2.93%  "
Set-TextValue "D22" "253.44"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("E24").Value = "  +6.96%  "
$ws.Range("E25").Value = "  -1.09%  "
Set-TextValue "D26" "9.82"
$ws.Range("E26").Value = "  +1.20%  "
Set-TextValue "D27" "171.21"
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("E28").Value = "  -2.30%  "
Set-TextValue "D29" "20.59"
$ws.Range("E29").Value = "  +2.46%  "
Set-TextValue "D30" "1.42"
$ws.Range("E30").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -0.18%  "
Set-TextValue "D32" "0.123"
$ws.Range("E32").Value = "  +0.20%  "
Set-TextValue "D33" "5.05"
$ws.Range("E33").Value = "  +0.44%  "
Set-TextValue "D34" "4.72"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("E35").Value = "  +2.94%  "
Set-TextValue "D36" "6.49"
$ws.Range("E36").Value = "  -3.11%  "
$ws.Range("E37").Value = "  +0.91%  "
Set-TextValue "D38" "3.59"
$ws.Range("E38").Value = "  -3.46%  "
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("E40").Value = "  -0.40%  "
Set-TextValue "D41" "8.74"
$ws.Range("E41").Value = "  +1.44%  "
Set-TextValue "D42" "0.000222"
$ws.Range("E42").Value = "  -13.78%  "
$ws.Range("B43").Value = "FTXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D43" "4.49"
$ws.Range("E43").Value = "  -7.22%  "
Set-TextValue "D44" "0.0964"
$ws.Range("E44").Value = "  -0.90%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D45" "1.21"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "98.19"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "1.476.69"
$ws.Range("E47").Value = "  -0.30%  "
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("E51").Value = "  +5.96%  "
